$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feb")

# Update data values (row numbers match the worksheet rows)
$ws.Range("C5").Value = 1

$ws.Range("B6").Value = 170
$ws.Range("C6").Value = 14
$ws.Range("F6").Value = 17.65

$ws.Range("B7").Value = 99
$ws.Range("F7").Value = 20.2

$ws.Range("B9").Value = 55
$ws.Range("F9").Value = 18.18

$ws.Range("B10").Value = 75

$ws.Range("B12").Value = 80

$ws.Range("B13").Value = 119

$ws.Range("B14").Value = 106

$ws.Range("B15").Value = 124

$ws.Range("B17").Value = 44
$ws.Range("F17").Value = 22.73

$ws.Range("B19").Value = 83

# Move the active selection to B10, matching the saved view state
$ws.Activate()
$ws.Range("B10").Select()
